$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.103.96"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "2.928.45"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.03"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.81%  "
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "3.413.66"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "61.104.02"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").Value = "2.928.39"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "434.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.27%  "
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("E31").Value = "  +3.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").Value = "0.0₃0867"
$ws.Range("E34").Value = "  +2.23%  "
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.284"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "374.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0346"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.07%  "
$ws.Range("D45").Value = "2.691.04"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.60%  "
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.105"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("E51").Value = "  -0.14%  "
